# Update NATMI TPM-derived Ligand/Receptor expression metrics (Vegfb -> Nrp1)
# with newly recomputed TPM values.
#
# Root cause of the change: re-running the scoring script with new TPM input
# changed two raw per-cluster averages:
#   - Ligand average expression (col G) for sending cluster "ECs"
#   - Receptor average expression (col M) for target cluster "ECs"
# Every other touched column (H, I, J, N, O, P, Q, R, S, T) is a downstream
# statistic derived from those two numbers (total = avg * expressing-cell
# count, specificity = avg/total divided by the cross-cluster sum, and edge
# weight/specificity = ligand-side * receptor-side). The literal values below
# are those recomputed downstream figures, applied directly to each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ G=1.038291333333333; H=3.114874;    I=0.09907705749147798; J=0.09907705749147799; M=123.2806423333333; N=369.841927; O=0.6241574062367528; P=0.6241574062367526; Q=128.0012225024664;  R=1152.011002522198; S=0.06183967922145053; T=0.06183967922145053 }
    3  = @{ G=1.038291333333333; H=3.114874;    I=0.09907705749147798; J=0.09907705749147799;                                    O=0.2392728888301323; P=0.2392728888301322; Q=49.06970898032777;  R=441.6273808229499; S=0.02370645376277503; T=0.02370645376277503 }
    4  = @{ G=1.038291333333333; H=3.114874;    I=0.09907705749147798; J=0.09907705749147799;                                    O=0.136569704933115;  P=0.136569704933115;  Q=28.00750101426978;  R=252.067509128428;  S=0.01353092450725242; T=0.01353092450725242 }
    5  = @{                                     I=0.4810466828311408;  J=0.4810466828311408;  M=123.2806423333333; N=369.841927; O=0.6241574062367528; P=0.6241574062367526; Q=621.4815522598511;  R=5593.33397033866;  S=0.3002488498346787;  T=0.3002488498346786 }
    6  = @{                                     I=0.4810466828311408;  J=0.4810466828311408;                                     O=0.2392728888301323; P=0.2392728888301322;                                             S=0.1151014294631594;  T=0.1151014294631594 }
    7  = @{                                     I=0.4810466828311408;  J=0.4810466828311408;                                     O=0.136569704933115;  P=0.136569704933115;                                              S=0.06569640353330268; T=0.06569640353330267 }
    8  = @{                                     I=0.4198762596773812;  J=0.4198762596773812;  M=123.2806423333333; N=369.841927; O=0.6241574062367528; P=0.6241574062367526; Q=542.4532772694693;  R=4882.079495425222; S=0.2620688771806235;  T=0.2620688771806234 }
    9  = @{                                     I=0.4198762596773812;  J=0.4198762596773812;                                     O=0.2392728888301323; P=0.2392728888301322;                                             S=0.1004650056041978;  T=0.1004650056041978 }
    10 = @{                                     I=0.4198762596773812;  J=0.4198762596773812;                                     O=0.136569704933115;  P=0.136569704933115;                                              S=0.05734237689255994; T=0.05734237689255993 }
}

$colNum = @{ G=7; H=8; I=9; J=10; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($r in $values.Keys) {
    $rowVals = $values[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Cells.Item([int]$r, $colNum[$col]).Value2 = $rowVals[$col]
    }
}
